$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for row-level fields (D, L, M, N, O, P, Q, R, S, T) across rows 3-13.
# The underlying edit is a permutation of whole "data rows" (while A,B,C,E,F,G,H,I,J,K
# stay identical for every row in this block), so we simply overwrite each row's
# changed columns with its final values.

$rows = @{
    3  = @{ D = 44572; L = 'Primera';  M = 65;  N = 20000; O = 20000; P = 20000; Q = '$/bandeja 6 kilos'; R = 'Región Metropolitana';                     S = 3333; T = 6 }
    4  = @{ D = 44187; L = 'Especial'; M = 45;  N = 14000; O = 14000; P = 14000; Q = '$/bandeja 7 kilos'; R = 'Provincia de San Felipe de Aconcagua';     S = 2000; T = 7 }
    5  = @{ D = 44187; L = 'Primera';  M = 50;  N = 12000; O = 12000; P = 12000; Q = '$/bandeja 7 kilos'; R = 'Provincia de San Felipe de Aconcagua';     S = 1714; T = 7 }
    6  = @{ D = 44189; L = 'Especial'; M = 20;  N = 15000; O = 15000; P = 15000; Q = '$/bandeja 7 kilos'; R = 'Provincia de San Felipe de Aconcagua';     S = 2143; T = 7 }
    7  = @{ D = 44189; L = 'Primera';  M = 30;  N = 13000; O = 13000; P = 13000; Q = '$/bandeja 7 kilos'; R = 'Provincia de San Felipe de Aconcagua';     S = 1857; T = 7 }
    8  = @{ D = 44553; L = 'Especial'; M = 200; N = 22000; O = 22000; P = 22000; Q = '$/bandeja 6 kilos'; R = 'Provincia de San Felipe de Aconcagua';     S = 3667; T = 6 }
    9  = @{ D = 44553; L = 'Primera';  M = 150; N = 18000; O = 18000; P = 18000; Q = '$/bandeja 6 kilos'; R = 'Provincia de San Felipe de Aconcagua';     S = 3000; T = 6 }
    10 = @{ D = 44204; L = 'Primera';  M = 110; N = 7000;  O = 7500;  P = 7318;  Q = '$/bandeja 7 kilos'; R = 'Provincia de San Felipe de Aconcagua';     S = 1045; T = 7 }
    11 = @{ D = 44550; L = 'Primera';  M = 60;  N = 24000; O = 24000; P = 24000; Q = '$/bandeja 7 kilos'; R = 'Región Metropolitana';                     S = 3429; T = 7 }
    12 = @{ D = 44558; L = 'Especial'; M = 20;  N = 22000; O = 22000; P = 22000; Q = '$/bandeja 6 kilos'; R = 'Provincia de San Felipe de Aconcagua';     S = 3667; T = 6 }
    13 = @{ D = 44558; L = 'Primera';  M = 25;  N = 18000; O = 18000; P = 18000; Q = '$/bandeja 6 kilos'; R = 'Provincia de San Felipe de Aconcagua';     S = 3000; T = 6 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("Q$r").Value = $vals.Q
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
    $ws.Range("T$r").Value = $vals.T
}
